$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B65").Value = "SingleUseId98"
$ws.Range("C65").Value = "Small"
$ws.Range("D65").Value = "Left"
$ws.Range("E65").Value = "LTR"
$ws.Range("F65").Value = "Fix: <value>"

$ws.Range("B66").Value = "SingleUseId99"
$ws.Range("C66").Value = "Small"
$ws.Range("D66").Value = "Left"
$ws.Range("E66").Value = "LTR"
$ws.Range("F66").Value = "NO FIX"
